# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to "_FV2210" / "_FV2304"
# - Turn the sheet's data range into an Excel Table ("Table1") with an AutoFilter
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header labels (without the _old/_new suffix) in column order
$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$colsFirst  = @("A","B","C","D","E","F","G","H","I","J")
$colsSecond = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Range($colsFirst[$i]  + "1").Value = $baseHeaders[$i] + "_FV2210"
    $ws.Range($colsSecond[$i] + "1").Value = $baseHeaders[$i] + "_FV2304"
}
# Column K1 keeps its existing "diff" label - no change needed.

# Convert the used range (including the freshly renamed header row) into a Table
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1) like the workbook now does
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
